# Insert a new data record as row 103, pushing the existing rows
# 103-147 down to 104-148 (dimension grows from A1:T147 to A1:T148).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(103).Insert()

$ws.Range("A103").Value = 11
$ws.Range("B103").Value = "Vega Monumental Concepción"
$ws.Range("C103").Value = "Bíobío"
$ws.Range("D103").Value = "2023-03-22"
$ws.Range("E103").Value = 8
$ws.Range("F103").Value = "Fruta"
$ws.Range("G103").Value = 100101
$ws.Range("H103").Value = "Berries"
$ws.Range("I103").Value = 100101001
$ws.Range("J103").Value = "Arándano (blue)"
$ws.Range("K103").Value = "Sin especificar"
$ws.Range("L103").Value = "Primera"
$ws.Range("M103").Value = 110
$ws.Range("N103").Value = 4000
$ws.Range("O103").Value = 4500
$ws.Range("P103").Value = 4227
$ws.Range("Q103").Value = "$/bandeja 2 kilos"
$ws.Range("R103").Value = "Provincia de Curicó"
$ws.Range("S103").Value = 2114
$ws.Range("T103").Value = 2
